$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F
$ws.Range("F1").Value = "Comment"

# Add the comment text for the Chinstrap/Adelie/Gentoo penguins row
$ws.Range("F8").Value = "Estimates from this study were originally calculated for the reproduction period only. They were raised to annual estimates using a simple cross product."
